$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.922.17"

$ws.Range("E2").Value = "  -2.40%  "

$ws.Range("D3").Value = "1.795.68"

$ws.Range("E3").Value = "  -0.40%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"

$ws.Range("E4").Value = "  +0.11%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "317.03"
$c.Style = "Normal"

$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("E6").Value = "  -0.02%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5317"
$c.Style = "Normal"

$ws.Range("E7").Value = "  -2.58%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3882"
$c.Style = "Normal"

$ws.Range("E8").Value = "  +2.90%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07440"
$c.Style = "Normal"

$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("E10").Value = "  -2.16%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.086"
$c.Style = "Normal"

$ws.Range("E11").Value = "  -2.56%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"

$ws.Range("E12").Value = "  +0.09%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.172"
$c.Style = "Normal"

$ws.Range("E13").Value = "  +0.37%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.431"
$c.Style = "Normal"

$ws.Range("E14").Value = "  +0.55%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "20.34"
$c.Style = "Normal"

$ws.Range("E15").Value = "  -1.51%  "

$ws.Range("D16").Value = "1.796.91"

$ws.Range("E16").Value = "  +0.10%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "88.35"
$c.Style = "Normal"

$ws.Range("E17").Value = "  -2.03%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001059"
$c.Style = "Normal"

$ws.Range("E18").Value = "  -0.60%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06540"
$c.Style = "Normal"

$ws.Range("E19").Value = "  +1.33%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"

$ws.Range("E20").Value = "  -0.03%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.23"
$c.Style = "Normal"

$ws.Range("E21").Value = "  -0.11%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.945"
$c.Style = "Normal"

$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("D23").Value = "27.963.02"

$ws.Range("E23").Value = "  -2.28%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.09"
$c.Style = "Normal"

$ws.Range("E24").Value = "  -0.22%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.090"
$c.Style = "Normal"

$ws.Range("E25").Value = "  -0.06%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "157.37"
$c.Style = "Normal"

$ws.Range("E26").Value = "  -0.88%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.12"
$c.Style = "Normal"

$ws.Range("E27").Value = "  -1.64%  "

$ws.Range("D28").Value = "2.001.03"

$ws.Range("E28").Value = "  -0.15%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.291"
$c.Style = "Normal"

$ws.Range("E29").Value = "  -2.58%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "121.79"
$c.Style = "Normal"

$ws.Range("E30").Value = "  -0.92%  "

$ws.Range("E31").Value = "  +2.56%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.097"
$c.Style = "Normal"

$ws.Range("E32").Value = "  -0.87%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.668"
$c.Style = "Normal"

$ws.Range("E33").Value = "  -0.35%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.491"
$c.Style = "Normal"

$ws.Range("E34").Value = "  -2.73%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.07031"
$c.Style = "Normal"

$ws.Range("E35").Value = "  +7.85%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.2202"
$c.Style = "Normal"

$ws.Range("E36").Value = "  -2.10%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02276"
$c.Style = "Normal"

$ws.Range("E37").Value = "  -1.25%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "5.075"
$c.Style = "Normal"

$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("E39").Value = "  -4.23%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "11.20"
$c.Style = "Normal"

$ws.Range("E40").Value = "  -0.49%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.183"
$c.Style = "Normal"

$ws.Range("E41").Value = "  -1.63%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.6107"
$c.Style = "Normal"

$ws.Range("E42").Value = "  -2.09%  "

$ws.Range("E43").Value = "  -0.89%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "13.37"
$c.Style = "Normal"

$ws.Range("E44").Value = "  +1.23%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.676"
$c.Style = "Normal"

$ws.Range("E45").Value = "  -0.44%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5699"
$c.Style = "Normal"

$ws.Range("E46").Value = "  -2.59%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "124.79"
$c.Style = "Normal"

$ws.Range("E47").Value = "  -1.40%  "

# Row 48 becomes NEARProtocol data
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.914"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "

# Row 49 becomes EOS data
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.177"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.74%  "

# Row 50
$ws.Range("E50").Value = "  -1.16%  "

# Row 51
$ws.Range("E51").Value = "  +30.49%  "
